$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format specific Price cells as Text so numeric-looking strings
# (e.g. "1.001", "314.69") are preserved verbatim instead of becoming numbers.
$textRows = @(5,7,8,9,10,11,12,13,15,16,17,18,19,20,22,25,26,27,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Apply the updated coin data
$ws.Cells.Item(2, 4).Value = '28.208.32'
$ws.Cells.Item(2, 5).Value = '  -0.17%  '
$ws.Cells.Item(3, 4).Value = '1.912.20'
$ws.Cells.Item(3, 5).Value = '  +2.03%  '
$ws.Cells.Item(4, 5).Value = '  -0.22%  '
$ws.Cells.Item(5, 4).Value = '314.69'
$ws.Cells.Item(5, 5).Value = '  +0.80%  '
$ws.Cells.Item(6, 5).Value = '  -0.13%  '
$ws.Cells.Item(7, 4).Value = '0.5081'
$ws.Cells.Item(7, 5).Value = '  +0.34%  '
$ws.Cells.Item(8, 4).Value = '0.3924'
$ws.Cells.Item(8, 5).Value = '  -0.46%  '
$ws.Cells.Item(9, 4).Value = '0.09303'
$ws.Cells.Item(9, 5).Value = '  -2.94%  '
$ws.Cells.Item(10, 4).Value = '1.141'
$ws.Cells.Item(10, 5).Value = '  -0.33%  '
$ws.Cells.Item(11, 4).Value = '41.87'
$ws.Cells.Item(11, 5).Value = '  +2.55%  '
$ws.Cells.Item(12, 4).Value = '6.390'
$ws.Cells.Item(12, 5).Value = '  -1.51%  '
$ws.Cells.Item(13, 4).Value = '20.92'
$ws.Cells.Item(13, 5).Value = '  -0.45%  '
$ws.Cells.Item(14, 4).Value = '1.904.17'
$ws.Cells.Item(14, 5).Value = '  +1.15%  '
$ws.Cells.Item(15, 4).Value = '7.313'
$ws.Cells.Item(15, 5).Value = '  -1.54%  '
$ws.Cells.Item(16, 4).Value = '1.001'
$ws.Cells.Item(16, 5).Value = '  -0.25%  '
$ws.Cells.Item(17, 4).Value = '0.00001121'
$ws.Cells.Item(17, 5).Value = '  -1.00%  '
$ws.Cells.Item(18, 4).Value = '92.45'
$ws.Cells.Item(18, 5).Value = '  -0.55%  '
$ws.Cells.Item(19, 4).Value = '0.06603'
$ws.Cells.Item(19, 5).Value = '  +0.02%  '
$ws.Cells.Item(20, 4).Value = '17.97'
$ws.Cells.Item(20, 5).Value = '  +2.05%  '
$ws.Cells.Item(22, 4).Value = '6.227'
$ws.Cells.Item(23, 4).Value = '28.268.17'
$ws.Cells.Item(23, 5).Value = '  -0.19%  '
$ws.Cells.Item(24, 5).Value = '  +1.16%  '
$ws.Cells.Item(25, 4).Value = '2.318'
$ws.Cells.Item(25, 5).Value = '  +0.60%  '
$ws.Cells.Item(26, 2).Value = 'LEO'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(26, 4).Value = '3.397'
$ws.Cells.Item(26, 5).Value = '  +0.22%  '
$ws.Cells.Item(27, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(27, 4).Value = '2.586'
$ws.Cells.Item(27, 5).Value = '  +0.65%  '
$ws.Cells.Item(28, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(28, 4).Value = '2.127.35'
$ws.Cells.Item(28, 5).Value = '  +1.45%  '
$ws.Cells.Item(29, 2).Value = 'EthereumClassic'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(29, 4).Value = '21.08'
$ws.Cells.Item(29, 5).Value = '  -0.65%  '
$ws.Cells.Item(30, 2).Value = 'Monero'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(30, 4).Value = '158.00'
$ws.Cells.Item(30, 5).Value = '  -0.54%  '
$ws.Cells.Item(31, 2).Value = 'BitcoinCash'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(31, 4).Value = '127.11'
$ws.Cells.Item(31, 5).Value = '  -0.39%  '
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).Value = '1.100'
$ws.Cells.Item(32, 5).Value = '  +3.06%  '
$ws.Cells.Item(33, 2).Value = 'Stellar'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(33, 4).Value = '0.1076'
$ws.Cells.Item(33, 5).Value = '  +0.56%  '
$ws.Cells.Item(34, 2).Value = 'Filecoin'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(34, 4).Value = '5.637'
$ws.Cells.Item(34, 5).Value = '  -0.09%  '
$ws.Cells.Item(35, 2).Value = 'HuobiToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(35, 4).Value = '3.616'
$ws.Cells.Item(35, 5).Value = '  -0.26%  '
$ws.Cells.Item(36, 2).Value = 'FraxShare'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(36, 4).Value = '9.714'
$ws.Cells.Item(36, 5).Value = '  +1.62%  '
$ws.Cells.Item(37, 2).Value = 'Hedera'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(37, 4).Value = '0.06673'
$ws.Cells.Item(37, 5).Value = '  -0.64%  '
$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(38, 4).Value = '0.02422'
$ws.Cells.Item(38, 5).Value = '  +1.25%  '
$ws.Cells.Item(39, 2).Value = 'ARBITRUM'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(39, 4).Value = '1.250'
$ws.Cells.Item(39, 5).Value = '  -0.82%  '
$ws.Cells.Item(40, 4).Value = '1.304'
$ws.Cells.Item(40, 5).Value = '  +9.76%  '
$ws.Cells.Item(41, 2).Value = 'Algorand'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(41, 4).Value = '0.2199'
$ws.Cells.Item(41, 5).Value = '  +0.30%  '
$ws.Cells.Item(42, 2).Value = 'TheSandbox'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(42, 4).Value = '0.6445'
$ws.Cells.Item(42, 5).Value = '  +1.02%  '
$ws.Cells.Item(43, 2).Value = 'Aptos'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(43, 4).Value = '11.48'
$ws.Cells.Item(43, 5).Value = '  -0.21%  '
$ws.Cells.Item(44, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(44, 4).Value = '4.989'
$ws.Cells.Item(44, 5).Value = '  -0.30%  '
$ws.Cells.Item(45, 2).Value = 'Frax'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(45, 4).Value = '1.000'
$ws.Cells.Item(45, 5).Value = '  -0.13%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).Value = '13.28'
$ws.Cells.Item(46, 5).Value = '  -1.99%  '
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).Value = '0.6033'
$ws.Cells.Item(47, 5).Value = '  +0.68%  '
$ws.Cells.Item(48, 2).Value = 'PancakeSwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(48, 4).Value = '3.721'
$ws.Cells.Item(48, 5).Value = '  +1.64%  '
$ws.Cells.Item(49, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(49, 4).Value = '1.283'
$ws.Cells.Item(49, 5).Value = '  +1.19%  '
$ws.Cells.Item(50, 2).Value = 'NEARProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(50, 4).Value = '2.021'
$ws.Cells.Item(50, 5).Value = '  +0.85%  '
$ws.Cells.Item(51, 2).Value = 'Quant'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(51, 4).Value = '122.86'
$ws.Cells.Item(51, 5).Value = '  -1.06%  '
